# Excel database connected with multiple screens.
# Now it can display all parameters based on the url.
#
# 1) Rename the original sheet to "Financials" and shift its existing
#    table down by 3 rows (row 3 -> row 6, rows 5-13 -> rows 8-16) by
#    inserting 3 blank rows above the header.
# 2) Add a second sheet "Demand_view" after "Financials" containing a
#    similarly-shaped table (shifted to columns D/F/H/J, rows 7 & 9-17)
#    with its own set of figures, and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet 1: "Financials" ------------------------------------------------
$ws1.Name = "Financials"

# Push the existing table down by 3 rows (3->6, 5->8 ... 13->16)
$ws1.Rows("3:5").Insert()

# Restore the selection to the (now relocated) table
$ws1.Range("C6:I16").Select()

# --- Sheet 2: "Demand_view" ------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Demand_view"

$ws2.Range("D7").Value = "id"
$ws2.Range("F7").Value = "label"
$ws2.Range("H7").Value = "percentage"
$ws2.Range("J7").Value = "title"

$ws2.Range("D9").Value = 1
$ws2.Range("F9").Value = "Budget"
$ws2.Range("H9").Value = 3
$ws2.Range("J9").Value = "Revenue"

$ws2.Range("D10").Value = 2
$ws2.Range("F10").Value = "Actuals"
$ws2.Range("H10").Value = 2
$ws2.Range("J10").Value = "Revenue"

$ws2.Range("D11").Value = 3
$ws2.Range("F11").Value = "Achieved percentage"
$ws2.Range("H11").Value = 44
$ws2.Range("J11").Value = "Revenue"

$ws2.Range("D12").Value = 4
$ws2.Range("F12").Value = "Budget"
$ws2.Range("H12").Value = 56
$ws2.Range("J12").Value = "Gross Margin $"

$ws2.Range("D13").Value = 5
$ws2.Range("F13").Value = "Actuals"
$ws2.Range("H13").Value = 88
$ws2.Range("J13").Value = "Gross Margin $"

$ws2.Range("D14").Value = 6
$ws2.Range("F14").Value = "Achieved percentage"
$ws2.Range("H14").Value = 8
$ws2.Range("J14").Value = "Gross Margin $"

$ws2.Range("D15").Value = 7
$ws2.Range("F15").Value = "Budget"
$ws2.Range("H15").Value = 7
$ws2.Range("J15").Value = "Gross Margin %"

$ws2.Range("D16").Value = 8
$ws2.Range("F16").Value = "Actuals"
$ws2.Range("H16").Value = 9
$ws2.Range("J16").Value = "Gross Margin %"

$ws2.Range("D17").Value = 9
$ws2.Range("F17").Value = "Achieved percentage"
$ws2.Range("H17").Value = 98
$ws2.Range("J17").Value = "Gross Margin %"

# Select F12 and make Demand_view the active/visible tab
$ws2.Range("F12").Select()
$ws2.Activate()
